$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the first occurrence of $needle inside the explicit range
# [$start, $end) with $replacement, WITHOUT wrapping past the range (Wrap=0 /
# wdFindStop) so the search cannot leak into neighbouring cells/paragraphs.
# ---------------------------------------------------------------------------
function Replace-InRange($start, $end, $needle, $replacement) {
  $rng = $d.Range($start, $end)
  $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 0, $false, $replacement, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) "Volume" / "Down" table cell -> collapse the two paragraphs into the
#    single (now empty) paragraph that used to hold "Volume".
# ---------------------------------------------------------------------------
$pVolume = $d.Paragraphs.Item(82)
$pDown   = $d.Paragraphs.Item(83)
Replace-InRange $pVolume.Range.Start $pDown.Range.End "Volume" ""
Replace-InRange $pVolume.Range.Start $pDown.Range.End "Down" ""
$d.Paragraphs.Item(83).Range.Delete()

# ---------------------------------------------------------------------------
# 2) "Volume" / "Up" table cell -> collapse the two paragraphs into the
#    single (now empty) paragraph that used to hold "Volume".
#    (indices shifted down by one paragraph after step 1's merge)
# ---------------------------------------------------------------------------
$pVolume2 = $d.Paragraphs.Item(83)
$pUp      = $d.Paragraphs.Item(84)
Replace-InRange $pVolume2.Range.Start $pUp.Range.End "Volume" ""
Replace-InRange $pVolume2.Range.Start $pUp.Range.End "Up" ""
$d.Paragraphs.Item(84).Range.Delete()

# ---------------------------------------------------------------------------
# 3) "Effects <tabs> Volume Control" paragraph -> strip every run after
#    "Effects" (the run of tab characters plus the trailing
#    "   Volume Control" text), leaving only the "Effects" run behind.
# ---------------------------------------------------------------------------
$effectsPara = $d.Paragraphs.Item(85)
$paraStart = $effectsPara.Range.Start
$paraEnd   = $effectsPara.Range.End

$findRng = $d.Range($paraStart, $paraEnd)
$findRng.Find.Execute("Effects", $false, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$afterEffects = $findRng.End

$tail = $d.Range($afterEffects, $paraEnd - 1)
$tail.Text = ""
